$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LEU")

# Row 47 - Long Term Investments
$ws.Range("D47").Value = 19700
$ws.Range("E47").Value = 29500
$ws.Range("F47").Value = 29800
$ws.Range("G47").Value = 34800
$ws.Range("H47").Value = 65600
$ws.Range("I47").Value = 48100
$ws.Range("J47").Value = 151300

# Row 52 - Other Assets
$ws.Range("D52").Value = 1100
$ws.Range("E52").Value = 24100
$ws.Range("F52").Value = 23600
$ws.Range("G52").Value = 183800
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 3600
$ws.Range("J52").Value = 24400

# Row 57 - Accounts Payable
$ws.Range("D57").Value = 48200

# Row 58 - Short/Current Long Term Debt
$ws.Range("D58").Value = 6100

# Row 83 - Depreciation (cash flow)
$ws.Range("J83").Value = "NA"

# Row 94 - Other Cashflows from Investing Activities
$ws.Range("J94").Value = "NA"

# Row 100 - Other Cash Flows from Financing Activities
$ws.Range("J100").Value = "NA"

# Row 101 - Effect Of Exchange Rate Changes
$ws.Range("J101").Value = "NA"
